# Updates cryptos list values (price & volume) per diff; B/C swap for rows 25-26.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.155.47"
$ws.Range("E2").Value = "'  +0.36%  "
$ws.Range("D3").Value = "'1.901.21"
$ws.Range("E3").Value = "'  +0.88%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("D5").Value = "'306.03"
$ws.Range("E5").Value = "'  -0.42%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "'  +0.17%  "
$ws.Range("D7").Value = "'0.5256"
$ws.Range("D8").Value = "'0.3774"
$ws.Range("E8").Value = "'  +1.35%  "
$ws.Range("D9").Value = "'0.07243"
$ws.Range("E9").Value = "'  +0.49%  "
$ws.Range("D10").Value = "'21.15"
$ws.Range("E10").Value = "'  +1.02%  "
$ws.Range("D11").Value = "'0.8975"
$ws.Range("E11").Value = "'  -0.66%  "
$ws.Range("D12").Value = "'0.08337"
$ws.Range("E12").Value = "'  +9.70%  "
$ws.Range("D13").Value = "'1.907.12"
$ws.Range("E13").Value = "'  +1.17%  "
$ws.Range("D14").Value = "'94.67"
$ws.Range("E14").Value = "'  -0.57%  "
$ws.Range("D15").Value = "'5.262"
$ws.Range("E15").Value = "'  -0.10%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "'  +0.09%  "
$ws.Range("D17").Value = "'0.000008623"
$ws.Range("E17").Value = "'  +1.40%  "
$ws.Range("D18").Value = "'14.52"
$ws.Range("D19").Value = "'0.9994"
$ws.Range("E19").Value = "'  +0.12%  "
$ws.Range("D20").Value = "'27.195.61"
$ws.Range("E20").Value = "'  +0.35%  "
$ws.Range("D21").Value = "'5.056"
$ws.Range("E21").Value = "'  +0.28%  "
$ws.Range("D22").Value = "'2.136.49"
$ws.Range("E22").Value = "'  +0.55%  "
$ws.Range("E23").Value = "'  +1.04%  "
$ws.Range("D24").Value = "'6.418"
$ws.Range("E24").Value = "'  -0.65%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'146.35"
$ws.Range("E25").Value = "'  +0.36%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.280"
$ws.Range("E26").Value = "'  +7.52%  "
$ws.Range("D27").Value = "'1.757"
$ws.Range("E27").Value = "'  -1.76%  "
$ws.Range("D28").Value = "'18.10"
$ws.Range("E28").Value = "'  +0.46%  "
$ws.Range("D29").Value = "'114.64"
$ws.Range("E29").Value = "'  -0.01%  "
$ws.Range("D30").Value = "'4.919"
$ws.Range("E30").Value = "'  -0.18%  "
$ws.Range("D31").Value = "'4.777"
$ws.Range("E31").Value = "'  -0.15%  "
$ws.Range("D32").Value = "'0.09262"
$ws.Range("E32").Value = "'  +0.55%  "
$ws.Range("D33").Value = "'0.8160"
$ws.Range("E33").Value = "'  +6.89%  "
$ws.Range("D34").Value = "'0.05047"
$ws.Range("E34").Value = "'  +0.09%  "
$ws.Range("D35").Value = "'1.234"
$ws.Range("E35").Value = "'  +3.90%  "
$ws.Range("D36").Value = "'2.976"
$ws.Range("E36").Value = "'  -1.27%  "
$ws.Range("D37").Value = "'3.338"
$ws.Range("E37").Value = "'  +2.02%  "
$ws.Range("D38").Value = "'2.589"
$ws.Range("E38").Value = "'  +2.57%  "
$ws.Range("D39").Value = "'0.5711"
$ws.Range("E39").Value = "'  +1.86%  "
$ws.Range("E40").Value = "'  -0.59%  "
$ws.Range("D41").Value = "'1.073"
$ws.Range("E41").Value = "'  -0.31%  "
$ws.Range("D42").Value = "'6.666"
$ws.Range("E42").Value = "'  +1.05%  "
$ws.Range("D43").Value = "'8.942"
$ws.Range("E43").Value = "'  +0.50%  "
$ws.Range("D44").Value = "'118.06"
$ws.Range("E44").Value = "'  -0.03%  "
$ws.Range("E45").Value = "'  +0.35%  "
$ws.Range("D46").Value = "'0.4832"
$ws.Range("E46").Value = "'  +0.80%  "
$ws.Range("D47").Value = "'0.9998"
$ws.Range("E47").Value = "'  +0.15%  "
$ws.Range("D48").Value = "'10.15"
$ws.Range("E48").Value = "'  +0.09%  "
$ws.Range("D49").Value = "'1.612"
$ws.Range("E49").Value = "'  +2.40%  "
$ws.Range("D50").Value = "'37.42"
$ws.Range("E50").Value = "'  +0.78%  "
$ws.Range("D51").Value = "'63.55"
$ws.Range("E51").Value = "'  -0.05%  "
